$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 1.62
$ws.Range("H2").Value = 3.9
$ws.Range("I2").Value = 5.5
$ws.Range("J2").Value = 2.25

$ws.Range("Q2").Value = 2.2
$ws.Range("R2").Value = 1.65
$ws.Range("S2").Value = 1.44
$ws.Range("T2").Value = 2.63

$ws.Range("Z2").Value = 11

$ws.Range("AC2").Value = 8.5
$ws.Range("AE2").Value = 23

$ws.Range("AK2").Value = 67
$ws.Range("AL2").Value = 51

$ws.Range("AO2").Value = 8.5

$ws.Range("AT2").Value = 2.63
